$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 278 (pushes old rows 278..336 down to 279..337,
# preserving all of their data/formatting) and populate it with the
# latest daily record for this market/product.
$ws.Rows(278).Insert()

$ws.Range("A278").Value = 10
$ws.Range("B278").Value = "Vega Modelo de Temuco"
$ws.Range("C278").Value = "La Araucanía"
$ws.Range("D278").Value = 44785
$ws.Range("E278").Value = 9
$ws.Range("F278").Value = 100114013
$ws.Range("G278").Value = "Zanahoria"
$ws.Range("H278").Value = "Sin especificar"
$ws.Range("I278").Value = "Primera"
$ws.Range("J278").Value = 100
$ws.Range("K278").Value = 9000
$ws.Range("L278").Value = 9000
$ws.Range("M278").Value = 9000
$ws.Range("N278").Value = "`$/saco 25 kilos"
$ws.Range("O278").Value = "Región de La Araucanía"
$ws.Range("P278").Value = 360
$ws.Range("Q278").Value = 25
$ws.Range("R278").Value = "Hortaliza"
